$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '67.370.51'
Set-TextValue "E2" '  -2.23%  '
Set-TextValue "D3" '3.252.07'
Set-TextValue "E3" '  -5.17%  '
Set-TextValue "D4" '0.998'
Set-TextValue "E4" '  -0.01%  '
Set-TextValue "D5" '586.68'
Set-TextValue "E5" '  -4.20%  '
Set-TextValue "D6" '147.11'
Set-TextValue "E6" '  -11.62%  '
Set-TextValue "D7" '0.999'
Set-TextValue "E7" '  -0.13%  '
Set-TextValue "D8" '3.246.47'
Set-TextValue "E8" '  -5.03%  '
Set-TextValue "D9" '0.537'
Set-TextValue "E9" '  -9.33%  '
Set-TextValue "E10" '  -13.13%  '
Set-TextValue "E11" '  -4.40%  '
Set-TextValue "D12" '0.498'
Set-TextValue "E12" '  -11.05%  '
Set-TextValue "D13" '0.0000243'
Set-TextValue "E13" '  -9.29%  '
Set-TextValue "D14" '37.72'
Set-TextValue "E14" '  -13.96%  '
Set-TextValue "D15" '3.769.95'
Set-TextValue "E15" '  -5.31%  '
Set-TextValue "D16" '67.351.70'
Set-TextValue "E16" '  -2.35%  '
Set-TextValue "D17" '3.250.32'
Set-TextValue "E17" '  -5.04%  '
Set-TextValue "E18" '  -6.05%  '
Set-TextValue "D19" '519.80'
Set-TextValue "E19" '  -9.86%  '
Set-TextValue "D20" '7.00'
Set-TextValue "E20" '  -13.44%  '
Set-TextValue "D21" '14.72'
Set-TextValue "E21" '  -13.61%  '
Set-TextValue "D22" '0.743'
Set-TextValue "E22" '  -11.60%  '
Set-TextValue "D23" '7.68'
Set-TextValue "E23" '  -13.57%  '
Set-TextValue "D24" '84.84'
Set-TextValue "E24" '  -10.99%  '
Set-TextValue "D25" '13.22'
Set-TextValue "E25" '  -11.98%  '
Set-TextValue "D26" '0.999'
Set-TextValue "E26" '  -0.21%  '
Set-TextValue "E27" '  -12.36%  '
Set-TextValue "E28" '  -11.79%  '
Set-TextValue "D29" '7.87'
Set-TextValue "E29" '  -7.73%  '
Set-TextValue "D30" '28.59'
Set-TextValue "E30" '  -12.09%  '
Set-TextValue "D31" '1.18'
Set-TextValue "E31" '  -4.21%  '
Set-TextValue "E32" '  -5.53%  '
Set-TextValue "D33" '6.45'
Set-TextValue "E33" '  -17.04%  '
Set-TextValue "B34" 'NEARProtocol'
Set-TextValue "C34" 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue "D34" '5.60'
Set-TextValue "E34" '  -14.30%  '
Set-TextValue "B35" 'FirstDigitalUSD'
Set-TextValue "C35" 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue "D35" '1.00'
Set-TextValue "E35" '  +0.09%  '
Set-TextValue "D36" '56.33'
Set-TextValue "E36" '  +0.94%  '
Set-TextValue "D37" '503.20'
Set-TextValue "E37" '  -13.98%  '
Set-TextValue "E38" '  -6.72%  '
Set-TextValue "D39" '0.0840'
Set-TextValue "E39" '  -11.53%  '
Set-TextValue "B40" 'Cosmos'
Set-TextValue "C40" 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue "D40" '8.79'
Set-TextValue "E40" '  -15.54%  '
Set-TextValue "B41" 'Kaspa'
Set-TextValue "C41" 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue "D41" '0.124'
Set-TextValue "E41" '  -11.76%  '
Set-TextValue "B42" 'dogwifhat'
Set-TextValue "C42" 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue "D42" '2.69'
Set-TextValue "E42" '  -12.92%  '
Set-TextValue "B43" 'Maker'
Set-TextValue "C43" 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue "D43" '2.901.48'
Set-TextValue "E43" '  -10.09%  '
Set-TextValue "D44" '0.262'
Set-TextValue "E44" '  -10.47%  '
Set-TextValue "E45" '  -8.78%  '
Set-TextValue "E46" '  -0.09%  '
Set-TextValue "D47" '26.20'
Set-TextValue "E47" '  -15.21%  '
Set-TextValue "E48" '  -17.29%  '
Set-TextValue "D49" '124.25'
Set-TextValue "E49" '  -6.18%  '
Set-TextValue "E50" '  -10.46%  '
Set-TextValue "E51" '  -18.13%  '
